$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.326.76'
$ws.Range("E2").Value = '  -5.55%  '

# Row 3
$ws.Range("D3").Value = '1.840.29'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  -0.43%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.87'
$ws.Range("E5").Value = '  +1.89%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  -0.26%  '

# Row 7
$ws.Range("E7").Value = '  -4.50%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3879'
$ws.Range("E8").Value = '  -5.31%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.97'
$ws.Range("E9").Value = '  -3.69%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07867'
$ws.Range("E10").Value = '  -4.32%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9695'
$ws.Range("E11").Value = '  -4.64%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.96'
$ws.Range("E12").Value = '  -8.56%  '

# Row 13
$ws.Range("D13").Value = '1.879.91'
$ws.Range("E13").Value = '  -4.38%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.791'
$ws.Range("E14").Value = '  -5.16%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.934'
$ws.Range("E15").Value = '  -5.48%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06904'
$ws.Range("E16").Value = '  +0.53%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.007'
$ws.Range("E17").Value = '  -0.38%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '87.42'

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000009998'
$ws.Range("E19").Value = '  -3.77%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.97'
$ws.Range("E20").Value = '  -4.81%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.004'
$ws.Range("E21").Value = '  -0.40%  '

# Row 22
$ws.Range("D22").Value = '28.389.50'
$ws.Range("E22").Value = '  -5.39%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.362'
$ws.Range("E23").Value = '  -5.81%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.13'
$ws.Range("E24").Value = '  -6.87%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.174'
$ws.Range("E25").Value = '  -1.21%  '

# Row 26
$ws.Range("D26").Value = '2.087.86'
$ws.Range("E26").Value = '  -5.06%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.63'
$ws.Range("E27").Value = '  -2.01%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.30'
$ws.Range("E28").Value = '  -4.00%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.903'
$ws.Range("E29").Value = '  -10.18%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.983'
$ws.Range("E30").Value = '  -6.30%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.14'
$ws.Range("E31").Value = '  -3.41%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9520'
$ws.Range("E32").Value = '  -6.89%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09357'
$ws.Range("E33").Value = '  -2.83%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.328'
$ws.Range("E34").Value = '  -5.30%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.463'
$ws.Range("E35").Value = '  -2.55%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.328'
$ws.Range("E36").Value = '  -6.62%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06046'
$ws.Range("E37").Value = '  -7.64%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02177'
$ws.Range("E38").Value = '  -5.37%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.156'
$ws.Range("E39").Value = '  -4.83%  '

# Row 40
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5663'
$ws.Range("E40").Value = '  -5.04%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.605'
$ws.Range("E41").Value = '  -4.47%  '

# Row 42
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.03'
$ws.Range("E42").Value = '  -6.62%  '

# Row 43
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1787'
$ws.Range("E43").Value = '  -3.60%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.393'
$ws.Range("E44").Value = '  -5.87%  '

# Row 45
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.225'
$ws.Range("E45").Value = '  -1.70%  '

# Row 46
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5332'
$ws.Range("E46").Value = '  -4.47%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.65'
$ws.Range("E47").Value = '  -6.83%  '

# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.07052'
$ws.Range("E48").Value = '  -6.66%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.859'
$ws.Range("E49").Value = '  -6.48%  '

# Row 50
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '112.74'
$ws.Range("E50").Value = '  -4.48%  '

# Row 51
$ws.Range("B51").Value = 'Chiliz'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GSCt2y6YSgO26+chiliz-chz'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1298'
$ws.Range("E51").Value = '  +0.27%  '
